$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = "dnasr281@gmail.com, System"
$replacement = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $target) {
        $cell.Value = $replacement
    }
}
